# DaySale report update:
#  - insert a new line item "CATAFLAM 75MG/3ML 6 AMP." as item #2 (pushing
#    IVYMOND SYRUP and TIRATAM down to #3 / #4)
#  - insert a new line item "سرنجات 3 سم" as item #5 at the end of the list
#    (before the totals / footer rows)
#  - update the totals cell and the footer timestamp accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel coercing a
# numeric-looking string ("120.00", "0:6", ...) into a number/time and
# without marking the cell with a quote-prefix (keeps the destination's
# existing style/number-format untouched, matching how the source report
# was generated).
function Set-TextValue {
    param($sheet, $rangeAddr, $text)
    $scratch = $sheet.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $sheet.Range($rangeAddr).PasteSpecial(-4163)
    $scratch.Clear()
}

# ---------------------------------------------------------------------
# 1) Insert the CATAFLAM row as the new row 8 (was: 7 BISOLOCK, 8 IVYMOND,
#    9 TIRATAM, 10 totals, 11 footer)
# ---------------------------------------------------------------------
$ws.Rows.Item(8).EntireRow.Insert()

# Recreate the merged cells for the new row (Insert does not carry them).
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

# Clone the row formatting from the row below (the shifted IVYMOND row).
$ws.Range("A9:Q9").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new item's data.
$ws.Range("A8").Value = 2
Set-TextValue $ws "C8" "CATAFLAM 75MG/3ML 6 AMP."
Set-TextValue $ws "H8" "0:6"
Set-TextValue $ws "L8" "1"
Set-TextValue $ws "N8" "120.00"
Set-TextValue $ws "P8" "19.2000"
Set-TextValue $ws "Q8" "0:1"

# Renumber the items that shifted down.
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4

# ---------------------------------------------------------------------
# 2) Insert the سرنجات 3 سم row as the new row 11 (currently holding the
#    totals row after step 1 shifted everything down once already).
# ---------------------------------------------------------------------
$ws.Rows.Item(11).EntireRow.Insert()

# Recreate the merged cells for the new row.
$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()

# Clone the row formatting from the row above (the TIRATAM row).
$ws.Range("A10:Q10").Copy()
$ws.Range("A11:Q11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new item's data.
$ws.Range("A11").Value = 5
Set-TextValue $ws "C11" "سرنجات 3 سم"
Set-TextValue $ws "H11" "0:0"
Set-TextValue $ws "L11" "0"
Set-TextValue $ws "N11" "2.00"
Set-TextValue $ws "P11" "2.0000"
Set-TextValue $ws "Q11" "1:0"

# ---------------------------------------------------------------------
# 3) Update totals and footer timestamp (rows shifted down to 12 / 13).
# ---------------------------------------------------------------------
$ws.Range("P12").Value = 313.2
Set-TextValue $ws "A13" "Wednesday, 10 September, 2025 9:35 AM"
